# Summer - Tuna Salad with apple and celery
# Reorder the "Dressing" and "Salad" sections on Sheet2 so that each
# section's header sits directly above its own ingredient rows:
#   rows 2-6   -> Dressing header + dressing ingredients
#   rows 7-11  -> Salad header + instructions
#   rows 12-17 -> Salad ingredients

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Wipe the existing content/formatting for the block we are rearranging so
# stray row heights / cell styles from the old layout do not linger.
$ws.Range("A2:H17").Clear()
$ws.Rows("2:17").AutoFit()

# --- Dressing section (was rows 7-13, now rows 2-6) -----------------------
$ws.Range("A2").Value = "Dressing"

$ws.Range("D3").Value = "sugar-free mayonnaise"
$ws.Range("E3").Value = "1 tbsp sugar-free mayonnaise"
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "tbsp"

$ws.Range("D4").Value = "plain Greek yogurt"
$ws.Range("E4").Value = "1 tbsp Greek yoghurt"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "tbsp"
$ws.Range("H4").Value = "DAIRY"

$ws.Range("D5").Value = "salt"
$ws.Range("E5").Value = "1/8 teaspoon salt (adjust to your preference)"
$ws.Range("F5").Value = 0.125
$ws.Range("G5").Value = "tsp"
$ws.Range("H5").Value = "SPICE"

$ws.Range("D6").Value = "pepper"
$ws.Range("E6").Value = "1/8 teaspoon pepper (adjust to your preference)"
$ws.Range("F6").Value = 0.125
$ws.Range("G6").Value = "tsp"
$ws.Range("H6").Value = "SPICE"

# --- Salad section (was rows 2-6, now rows 7-11) ---------------------------
$ws.Range("A7").Value = "Salad"

$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "Put the tuna in a bowl and break it up with a fork."

$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "Slice the celery and the apples to the same thickness. Chop the pickles small. "

$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "Add the finely chopped onion and the capers. `nMix it together with the dressing ingredients and done."
$ws.Rows("10:10").RowHeight = 30
$ws.Range("C10").WrapText = $true

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = "You can add any fresh herb you like. Dill goes especially well with tuna."

# --- Salad ingredients (was rows 14-17 & 7, now rows 12-17) ----------------
$ws.Range("D12").Value = "tuna in water"
$ws.Range("E12").Value = "2 cans of tuna in water, drained"
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = "can"

$ws.Range("D13").Value = "apple"
$ws.Range("E13").Value = "1/2 apple"
$ws.Range("F13").Value = 0.5
$ws.Range("G13").Value = "apple"
$ws.Range("H13").Value = "PRODUCE"

$ws.Range("D14").Value = "celery"
$ws.Range("E14").Value = "1 stalk celery"
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = "stalk"
$ws.Range("H14").Value = "PRODUCE"

$ws.Range("D15").Value = "capers"
$ws.Range("E15").Value = "1 tbsp capers"
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = "tbsp"

$ws.Range("D16").Value = "red onion"
$ws.Range("E16").Value = "1 spring onion or a small red onion, finely chopped"
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = "onion"
$ws.Range("H16").Value = "PRODUCE"

$ws.Range("D17").Value = "cornichon"
$ws.Range("E17").Value = "4 cornichon (small sugar free pickled cucumbers)"
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = "cornichon"

# Restore the selection shown in the saved workbook.
$ws.Range("C17").Select()
